$d = $word.ActiveDocument

# Anchor on the unique, unchanged text "气瓶类型；NO" which immediately precedes the
# block of NO / NO2 legend entries we need to edit (UI: NO2 -> NOx). After
# Find.Execute the range collapses to the end of the match, giving a stable
# base offset to compute the remaining character-level edits from.
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = "气瓶类型；NO"
$find.MatchWildcards = $false
$found = $find.Execute()
if (-not $found) {
    throw "anchor text not found"
}
$base = $rng.End

# Relative to $base the legend text reads (subscript shown as [2]):
#   浓度，NO[2]浓度；NO校正系数，NO[2]校正系数；NO修正值，NO[2]修正值；压力
#   0123 456 7                                              32 33343536
#
# index 5  -> first  "2" (subscript) in "NO2浓度"     -> becomes "X" (NOx concentration)
# index 18 -> second "2" (subscript) in "NO2校正系数"  -> becomes "X" (NOx correction factor)
# index 32 -> third  "2" (subscript) stays "2" (NO2 correction value kept), and a
#             new "，NOX修正值" clause (NOx correction value) is appended right
#             after it, i.e. inserted at offset 36, just before "；压力".

# 1) Insert the new "，NOX修正值" clause at offset 36, as plain (non-subscript)
#    text first so it merges cleanly into the neighbouring runs, then mark
#    only the "X" character subscript afterwards by absolute offset.
$checkIns = $d.Range($base + 36, $base + 37)
if ($checkIns.Text -ne "；") {
    throw "unexpected text at insertion point: " + $checkIns.Text
}

$insPoint = $d.Range($base + 36, $base + 36)
$insPoint.InsertAfter("，NO")

$insPoint2 = $d.Range($base + 39, $base + 39)
$insPoint2.InsertAfter("X")

$insPoint3 = $d.Range($base + 40, $base + 40)
$insPoint3.InsertAfter("修正值")

$xNew = $d.Range($base + 39, $base + 40)
if ($xNew.Text -ne "X") {
    throw "unexpected text for new NOx subscript: " + $xNew.Text
}
$xNew.Font.Subscript = -1

# 2) Second "2" (offset 18, in "NO2校正系数；") -> "X", keep subscript formatting.
$r2 = $d.Range($base + 18, $base + 19)
if ($r2.Text -ne "2") {
    throw "unexpected text at second NO2: " + $r2.Text
}
$r2.Text = "X"
$r2.Font.Subscript = -1

# 3) First "2" (offset 5, in "NO2浓度；") -> "X", keep subscript formatting.
$r1 = $d.Range($base + 5, $base + 6)
if ($r1.Text -ne "2") {
    throw "unexpected text at first NO2: " + $r1.Text
}
$r1.Text = "X"
$r1.Font.Subscript = -1

Write-Output "NO2 to NOx edit applied"
